# Applies the "Updated cryptos list" data refresh (Wed Jan 31 07:53:43 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.023.66"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "2.338.33"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("E5").Value = "  -1.44%  "

# Row 6
$ws.Range("D6").Formula = "'101.11"
$ws.Range("E6").Value = "  -1.79%  "

# Row 7
$ws.Range("D7").Formula = "'0.513"
$ws.Range("E7").Value = "  -4.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("D9").Formula = "'0.512"
$ws.Range("E9").Value = "  -3.46%  "

# Row 10
$ws.Range("D10").Formula = "'34.93"
$ws.Range("E10").Value = "  -2.27%  "

# Row 11
$ws.Range("D11").Formula = "'52.50"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12
$ws.Range("E12").Value = "  -1.98%  "

# Row 13
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("E14").Value = "  -2.42%  "

# Row 15
$ws.Range("D15").Formula = "'15.81"
$ws.Range("E15").Value = "  +5.29%  "

# Row 16
$ws.Range("D16").Value = "2.343.29"
$ws.Range("E16").Value = "  +1.26%  "

# Row 17
$ws.Range("E17").Value = "  +2.46%  "

# Row 18
$ws.Range("D18").Value = "42.952.02"
$ws.Range("E18").Value = "  -0.76%  "

# Row 19
$ws.Range("D19").Formula = "'6.23"
$ws.Range("E19").Value = "  +0.73%  "

# Row 20
$ws.Range("E20").Value = "  -4.45%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0909"
$ws.Range("E21").Value = "  -2.59%  "

# Row 22
$ws.Range("D22").Formula = "'68.06"
$ws.Range("E22").Value = "  -0.14%  "

# Row 23
$ws.Range("D23").Formula = "'236.63"
$ws.Range("E23").Value = "  -2.16%  "

# Row 24
$ws.Range("E24").Value = "  +0.35%  "

# Row 25
$ws.Range("E25").Value = "  -2.09%  "

# Row 26
$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").Formula = "'25.67"
$ws.Range("E27").Value = "  +3.06%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Formula = "'2.32"
$ws.Range("E28").Value = "  +1.04%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Formula = "'35.83"
$ws.Range("E29").Value = "  -2.52%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Formula = "'9.32"
$ws.Range("E30").Value = "  -3.53%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Formula = "'162.73"
$ws.Range("E31").Value = "  -4.95%  "

# Row 32
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Formula = "'1.00"
$ws.Range("E32").Value = "  -0.01%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Formula = "'5.12"
$ws.Range("E33").Value = "  -3.01%  "

# Row 34
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Formula = "'17.56"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35
$ws.Range("D35").Formula = "'4.64"
$ws.Range("E35").Value = "  +7.54%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Formula = "'0.0727"
$ws.Range("E36").Value = "  -2.12%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Formula = "'2.44"
$ws.Range("E37").Value = "  -4.10%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Formula = "'1.86"
$ws.Range("E38").Value = "  -1.78%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Formula = "'2.91"
$ws.Range("E39").Value = "  -4.95%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Formula = "'0.102"
$ws.Range("E40").Value = "  -3.01%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Formula = "'0.113"
$ws.Range("E41").Value = "  -2.65%  "

# Row 42
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Formula = "'2.57"
$ws.Range("E42").Value = "  +8.88%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.023.15"
$ws.Range("E43").Value = "  +2.24%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Formula = "'0.0285"
$ws.Range("E44").Value = "  -2.24%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "'18.94"
$ws.Range("E45").Value = "  -1.29%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Formula = "'10.16"
$ws.Range("E46").Value = "  +1.52%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Formula = "'2.94"
$ws.Range("E47").Value = "  -1.91%  "

# Row 48
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Formula = "'56.74"
$ws.Range("E48").Value = "  +2.39%  "

# Row 49
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Formula = "'2.90"
$ws.Range("E49").Value = "  -1.72%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.563.83"
$ws.Range("E50").Value = "  +0.99%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Formula = "'4.70"
$ws.Range("E51").Value = "  +1.34%  "

